# Seccion 2 del Egel
#
# This script applies three related edits to the document:
#   1) In the "print(funcion(5, 3))" code sample paragraph, the "print"/"("
#      and "funcion"/"(5, 3))" runs are merged into two runs ("print(" and
#      "funcion(5, 3))"), the spell-check proofErr markers around "print"
#      and "funcion" are removed (grammar markers are kept), and the
#      paragraph (plus the following blank paragraph and the "a) 5" answer
#      paragraph) get an explicit en-US language tag.
#   2) The correct answer to question 7 ("d) Arbol binario") is bolded.
#   3) The correct answer to question 13 ("b) ... se queda inactivo ...")
#      is bolded.

$d = $word.ActiveDocument

function Set-ParagraphOuterXml {
    # NOTE: called positionally -- passing COM objects via named (-Param)
    # binding does not forward them correctly in this PowerShell runtime.
    param($Paragraph, [string]$InnerParagraphXml)

    $packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $InnerParagraphXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $Paragraph.Range.InsertXML($packageXml)
}

# ---------------------------------------------------------------------------
# 1) "print(funcion(5, 3))" paragraph + the two following paragraphs.
# ---------------------------------------------------------------------------

$codeParagraph = $null
$codeIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*print(funcion(5, 3))*") {
        $codeParagraph = $candidate
        $codeIndex = $i
        break
    }
}

$printInner = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:rPr>' +
    '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/>' +
    '</w:rPr></w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr>' +
    '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/>' +
    '</w:rPr><w:t>print(</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr>' +
    '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/>' +
    '</w:rPr><w:t>funcion(5, 3))</w:t></w:r>' +
    '</w:p>'
Set-ParagraphOuterXml $codeParagraph $printInner

$blankParagraph = $d.Paragraphs.Item($codeIndex + 1)
$blankInner = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:rPr>' +
    '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/>' +
    '</w:rPr></w:pPr>' +
    '</w:p>'
Set-ParagraphOuterXml $blankParagraph $blankInner

$answerAParagraph = $d.Paragraphs.Item($codeIndex + 2)
$answerAInner = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:rPr>' +
    '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/>' +
    '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' +
    '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
    '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/>' +
    '</w:rPr><w:t>a) 5</w:t></w:r>' +
    '</w:p>'
Set-ParagraphOuterXml $answerAParagraph $answerAInner

# ---------------------------------------------------------------------------
# 2) and 3) Bold the correct answers for question 7 and question 13.
# ---------------------------------------------------------------------------

function Add-BoldToParagraph {
    # NOTE: called positionally, see comment above.
    param($Paragraph)

    $Paragraph.Range.Font.Bold = 1
    $Paragraph.Range.Font.BoldBi = 1
}

$arbolParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*rbol binario*") {
        $arbolParagraph = $candidate
        break
    }
}
Add-BoldToParagraph $arbolParagraph

$deadlockParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*queda inactivo*") {
        $deadlockParagraph = $candidate
        break
    }
}
Add-BoldToParagraph $deadlockParagraph

Write-Host "Edits applied"
